# Re-apply the "No Style, No Grid" table style to the three tables that
# previously used the custom "Table_0" style ({542C6F6E-60A9-4974-AE94-182D45706913}).
# This mirrors what happens in PowerPoint's UI when a user selects a table
# and picks a different built-in style from the Table Styles gallery.

$p = $ppt.ActivePresentation

$oldStyle = "{542C6F6E-60A9-4974-AE94-182D45706913}"
$newStyle = "{670CC75C-E6AD-4F19-8013-F6C31A58B0D9}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyle) {
                $tbl.ApplyStyle($newStyle)
            }
        }
    }
}
